$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look numeric stay as text, matching the
# original inline-string cell type, by forcing a Text number format before assignment.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.497.30"
$ws.Range("E2").Value = "  +2.58%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.144.47"
$ws.Range("E3").Value = "  +2.32%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.64"
$ws.Range("E5").Value = "  +0.46%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.16"
$ws.Range("E6").Value = "  +6.17%  "

# Row 7
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.143.06"
$ws.Range("E8").Value = "  +2.40%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").Value = "  +1.83%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.57"
$ws.Range("E10").Value = "  +3.26%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("E11").Value = "  +2.21%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.474"
$ws.Range("E12").Value = "  +1.48%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000243"
$ws.Range("E13").Value = "  +1.50%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.96"
$ws.Range("E14").Value = "  +3.30%  "

# Row 15
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "68.483.27"
$ws.Range("E15").Value = "  +2.66%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.671.24"
$ws.Range("E16").Value = "  +2.48%  "

# Row 17
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.122"
$ws.Range("E17").Value = "  +1.00%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.13"
$ws.Range("E18").Value = "  +1.93%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.140.35"
$ws.Range("E19").Value = "  +1.71%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.75"
$ws.Range("E20").Value = "  -1.02%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "490.79"
$ws.Range("E21").Value = "  -0.34%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.700"
$ws.Range("E22").Value = "  +1.57%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.82"
$ws.Range("E23").Value = "  +1.22%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.05"
$ws.Range("E24").Value = "  +1.46%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.03"
$ws.Range("E25").Value = "  +2.72%  "

# Row 26
$ws.Range("E26").Value = "  +5.86%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.55"
$ws.Range("E27").Value = "  +3.77%  "

# Row 28
$ws.Range("E28").Value = "  +0.07%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.13"
$ws.Range("E29").Value = "  +4.46%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.37"
$ws.Range("E30").Value = "  +4.09%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.64"
$ws.Range("E31").Value = "  +1.24%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.43"
$ws.Range("E32").Value = "  +2.88%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.114"
$ws.Range("E33").Value = "  +1.67%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0959"
$ws.Range("E34").Value = "  +5.33%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.10%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "49.38"
$ws.Range("E36").Value = "  +5.47%  "

# Row 37
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.67"
$ws.Range("E37").Value = "  +1.67%  "

# Row 38
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.959"
$ws.Range("E38").Value = "  +1.05%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.323"
$ws.Range("E39").Value = "  +7.26%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.07"
$ws.Range("E40").Value = "  +5.28%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "49.14"
$ws.Range("E41").Value = "  +0.13%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.124"
$ws.Range("E42").Value = "  +1.43%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.44"
$ws.Range("E43").Value = "  +1.48%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.71"
$ws.Range("E44").Value = "  +8.60%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "395.94"
$ws.Range("E45").Value = "  +7.88%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.783.48"
$ws.Range("E46").Value = "  +1.00%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.22"
$ws.Range("E47").Value = "  +10.49%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0350"
$ws.Range("E48").Value = "  +1.27%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "135.64"
$ws.Range("E49").Value = "  +0.19%  "

# Row 50
$ws.Range("E50").Value = "  +0.01%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.35"
$ws.Range("E51").Value = "  +8.93%  "
